# A 17788-2020.xlsx — automatic update of files.
#
# The edit swaps the full record content between row-pairs (4,7), (5,8) and
# (6,9) on the "Artfynd" sheet: what used to be rows 7-9 becomes rows 4-6,
# and what used to be rows 4-6 becomes rows 7-9 (each row keeps its own
# formatting/row number, only the species-record data moves).
#
# Rows 4-6 additionally lose their placeholder empty J/K/N/AF cells (they
# get cleared), because the data that replaces them never carried those
# placeholder cells; rows 7-9 conversely gain that placeholder pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value. A leading "'" forces text so values
# that look like dates (e.g. "2022-10-21") are not reinterpreted as date
# serial numbers.
function Set-Text {
    param($Sheet, [string]$Addr, [string]$Text)
    $Sheet.Range($Addr).Value2 = "'" + $Text
}

# Helper: write a literal numeric value.
function Set-Num {
    param($Sheet, [string]$Addr, $Number)
    $Sheet.Range($Addr).Value2 = $Number
}

# ---- Row 4 <- (old) Row 7 content ----------------------------------------
Set-Num  $ws "A4"  110200520
Set-Num  $ws "B4"  77259
Set-Num  $ws "E4"  228912
Set-Text $ws "F4"  "Mörk kolflarnlav"
Set-Text $ws "G4"  "Carbonicola myrmecina"
Set-Text $ws "H4"  "(Ach.) Bendiksby & Timdal"
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("N4").ClearContents()
Set-Text $ws "P4"  "Storhöjden, Ång"
Set-Num  $ws "Q4"  650679.6700574562
Set-Num  $ws "R4"  7006329.71806069
Set-Num  $ws "S4"  10
Set-Text $ws "Y4"  "2022-10-21"
Set-Text $ws "AA4" "2022-10-21"
$ws.Range("AF4").ClearContents()
Set-Text $ws "AW4" "klara linder"
Set-Text $ws "AX4" "klara linder, Emmy Ransgart, Astrid Hedman, Ulf Sperens, Nils Ericson"

# ---- Row 5 <- (old) Row 8 content ----------------------------------------
Set-Num  $ws "A5"  110204944
Set-Num  $ws "B5"  77506
Set-Text $ws "D5"  "NT"
Set-Num  $ws "E5"  6425
Set-Text $ws "F5"  "Garnlav"
Set-Text $ws "G5"  "Alectoria sarmentosa"
Set-Text $ws "H5"  "(Ach.) Ach."
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("N5").ClearContents()
Set-Text $ws "P5"  "Vitberget, Ång"
Set-Num  $ws "Q5"  650712.081957066
Set-Num  $ws "R5"  7006373.321625493
Set-Num  $ws "S5"  10
Set-Text $ws "Y5"  "2022-10-21"
Set-Text $ws "AA5" "2022-10-21"
$ws.Range("AF5").ClearContents()
Set-Text $ws "AW5" "klara linder"
Set-Text $ws "AX5" "klara linder, Emmy Ransgart, Astrid Hedman, Ulf Sperens, Jonas Forsberg, Nils Ericson"

# ---- Row 6 <- (old) Row 9 content ----------------------------------------
Set-Num  $ws "A6"  110200634
Set-Num  $ws "B6"  81236
Set-Num  $ws "E6"  1312
Set-Text $ws "F6"  "Gammelgransskål"
Set-Text $ws "G6"  "Pseudographis pinicola"
Set-Text $ws "H6"  "(Nyl.) Rehm"
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("N6").ClearContents()
Set-Text $ws "P6"  "Storhöjden, Ång"
Set-Num  $ws "Q6"  650713.7645427607
Set-Num  $ws "R6"  7006376.115701492
Set-Num  $ws "S6"  10
Set-Text $ws "Y6"  "2022-10-21"
Set-Text $ws "AA6" "2022-10-21"
$ws.Range("AF6").ClearContents()
Set-Text $ws "AW6" "klara linder"
Set-Text $ws "AX6" "klara linder, Emmy Ransgart, Astrid Hedman, Ulf Sperens, Jonas Forsberg, Nils Ericson"

# ---- Row 7 <- (old) Row 4 content ----------------------------------------
Set-Num  $ws "A7"  97272055
Set-Num  $ws "B7"  85703
Set-Num  $ws "E7"  510
Set-Text $ws "F7"  "Doftskinn"
Set-Text $ws "G7"  "Cystostereum murrayi"
Set-Text $ws "H7"  "(Berk. & M.A. Curtis.) Pouzar"
Set-Text $ws "P7"  "Brattfaret, Ång"
Set-Num  $ws "Q7"  651119.2294871274
Set-Num  $ws "R7"  7005953.279259816
Set-Num  $ws "S7"  5
Set-Text $ws "Y7"  "2021-10-26"
Set-Text $ws "AA7" "2021-10-26"
Set-Text $ws "AW7" "Johannes Esberg"
Set-Text $ws "AX7" "Johannes Esberg"

# ---- Row 8 <- (old) Row 5 content ----------------------------------------
Set-Num  $ws "A8"  97272062
Set-Num  $ws "B8"  78596
Set-Text $ws "D8"  "LC"
Set-Num  $ws "E8"  6462
Set-Text $ws "F8"  "Stuplav"
Set-Text $ws "G8"  "Nephroma bellum"
Set-Text $ws "H8"  "(Spreng.) Tuck."
Set-Text $ws "P8"  "Brattfaret, Ång"
Set-Num  $ws "Q8"  651134.4035808664
Set-Num  $ws "R8"  7005919.585986075
Set-Num  $ws "S8"  5
Set-Text $ws "Y8"  "2021-10-26"
Set-Text $ws "AA8" "2021-10-26"
Set-Text $ws "AW8" "Johannes Esberg"
Set-Text $ws "AX8" "Johannes Esberg"

# ---- Row 9 <- (old) Row 6 content ----------------------------------------
Set-Num  $ws "A9"  97272044
Set-Num  $ws "B9"  77506
Set-Num  $ws "E9"  6425
Set-Text $ws "F9"  "Garnlav"
Set-Text $ws "G9"  "Alectoria sarmentosa"
Set-Text $ws "H9"  "(Ach.) Ach."
Set-Text $ws "P9"  "Brattfaret, Ång"
Set-Num  $ws "Q9"  651092.6504968955
Set-Num  $ws "R9"  7006066.106331692
Set-Num  $ws "S9"  5
Set-Text $ws "Y9"  "2021-10-26"
Set-Text $ws "AA9" "2021-10-26"
Set-Text $ws "AW9" "Johannes Esberg"
Set-Text $ws "AX9" "Johannes Esberg"
